$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Fill in missing Number for existing row 31 (numeric-looking string -> must stay text)
Set-TextValue $ws.Range("F31") "54"

# New rows 32-41
$rows = @(
    @{ Row = 32; A = $null; C = "439_九星叶_undefined_undefined_1bunch"; F = "20" },
    @{ Row = 33; A = "12";  C = "12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "85" },
    @{ Row = 34; A = "13";  C = "300_白星_White Gypso_ gypsophila_1kg"; F = "56" },
    @{ Row = 35; A = "14";  C = "300_白星_White Gypso_ gypsophila_1kg"; F = "56" },
    @{ Row = 36; A = "15";  C = "300_白星_White Gypso_ gypsophila_1kg"; F = "23" },
    @{ Row = 37; A = $null; C = "1_白洋桔梗_White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "55" },
    @{ Row = 38; A = "16";  C = "522_山归来绿_Smilax china_undefined_1bunch"; F = "5" },
    @{ Row = 39; A = $null; C = "448_吊米 绿_hanging amaranthus`ngreen_undefined_1bunch"; F = "12" },
    @{ Row = 40; A = $null; C = "325_小盼草_Northern Sea Oats_undefined_1bunch"; F = "20" },
    @{ Row = 41; A = $null; C = "320_雪柳花_Spiraea flower white_undefined_1bunch"; F = $null }
)

foreach ($item in $rows) {
    if ($item.A -ne $null) {
        Set-TextValue $ws.Cells.Item($item.Row, 1) $item.A
    }
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    if ($item.F -ne $null) {
        Set-TextValue $ws.Cells.Item($item.Row, 6) $item.F
    }
}

# Update the Summary sheet's concatenated digit string in G2 (huge numeric string -> must stay text)
Set-TextValue $summary.Range("G2") "015196181942320232115225241410308117766324040401156054208556562355512200"
